$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark after "Government shutdown from "
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. In the table, row 2 (NY Times row):
#    - cell (2,4) "Programmed" (currently "Kavanaugh") gets a new paragraph "Shut down"
#    - cell (2,5) "Downloaded" (currently empty) becomes "BOTH" with the _GoBack bookmark
#      placed right after the text (collapsed, not wrapping it)
$t = $d.Tables.Item(1)

$progCell = $t.Cell(2, 4)
$progRange = $progCell.Range
$newline = [char]13
$progRange.InsertAfter($newline + "Shut down")

$dlCell = $t.Cell(2, 5)
$dlRange = $dlCell.Range
$dlRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $dlRange)
$dlRange.InsertBefore("BOTH")
